$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted at the top of the "Acelga" price
# history (row 177), which pushes every existing record (rows 177-211)
# down by one row; the last existing record (row 211) lands on a brand
# new row 212.

# Step 1: create the new last row (212) by duplicating the full row 211
# (all of its static columns - market, region, product, unit, origin,
# etc. - are identical for every record in this sheet).
$ws.Range("A211:R211").Copy($ws.Range("A212:R212"))

# Step 2: shift the per-record data (Fecha, Volumen, Precio minimo,
# Precio maximo, Precio promedio ponderado, Precio $/Kg) down by one
# row, working from the bottom up so we never overwrite a source row
# before it has been copied.
for ($r = 211; $r -ge 178; $r--) {
    $srcRow = $r - 1
    $ws.Range("D$srcRow").Copy($ws.Range("D$r"))
    $ws.Range("J${srcRow}:M$srcRow").Copy($ws.Range("J${r}:M$r"))
    $ws.Range("P$srcRow").Copy($ws.Range("P$r"))
}

# Step 3: write the brand-new weekly record into row 177.
$ws.Range("D177").Value = 44637
$ws.Range("J177").Value = 200
$ws.Range("K177:M177").Value = 4000
$ws.Range("P177").Value = 1000
